{"js": "// Collapse the exploded \"[coordinatorXXX]\" run sequence (separated by\n// proofErr spell-check wrappers and a line break) down to a single\n// \"[coordinator]\" placeholder, and likewise collapse the comma-joined\n// \"[contactXXX]\" placeholder list down to a single \"[contact]\" placeholder.\n// (Christiane review #73056.)\n\nconst body = context.document.body;\n\n// --- Project Coordinator / Principal Investigator cell -------------------\n// Original text (the \"\\u000b\" is the <w:br/> line break between the email\n// row and the affiliation row):\n//   [coordinatorname], [coordinatormail], [coordinatorid], \\u000b[coordinatoraffiliation], [coordinatorror]\nconst coordinatorResults = body.search(\n  \"[coordinatorname], [coordinatormail], [coordinatorid], \\u000b[coordinatoraffiliation], [coordinatorror]\",\n  { matchCase: true }\n);\ncoordinatorResults.load(\"items\");\nawait context.sync();\n\nif (coordinatorResults.items.length > 0) {\n  coordinatorResults.items[0].insertText(\"[coordinator]\", \"Replace\");\n  await context.sync();\n}\n\n// --- Contact person cell --------------------------------------------------\nconst contactResults = body.search(\n  \"[contactname], [contactmail], [contactid], [contactaffiliation], [contactror]\",\n  { matchCase: true }\n);\ncontactResults.load(\"items\");\nawait context.sync();\n\nif (contactResults.items.length > 0) {\n  contactResults.items[0].insertText(\"[contact]\", \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Collapse the exploded \"[coordinatorXXX]\" run sequence (separated by\n# proofErr spell-check wrappers and a line break) down to a single\n# \"[coordinator]\" placeholder, and likewise collapse the comma-joined\n# \"[contactXXX]\" placeholder list down to a single \"[contact]\" placeholder.\n# (Christiane review #73056.)\n\n$d = $word.ActiveDocument\n\n# --- Project Coordinator / Principal Investigator cell -------------------\n# The original run sequence reads as one continuous string once rendered,\n# with the mid-paragraph <w:br/> showing up as a vertical-tab (chr 11) in\n# Word's Find text:\n#   [coordinatorname], [coordinatormail], [coordinatorid], <br/>[coordinatoraffiliation], [coordinatorror]\n$coordinatorNeedle = \"[coordinatorname], [coordinatormail], [coordinatorid], \" + [char]11 + \"[coordinatoraffiliation], [coordinatorror]\"\n$coordinatorFind = $d.Content.Find\n$coordinatorFind.ClearFormatting()\n$coordinatorFind.Replacement.ClearFormatting()\n$coordinatorFind.Execute(\n    $coordinatorNeedle,   # FindText\n    $false,               # MatchCase\n    $false,               # MatchWholeWord\n    $false,               # MatchWildcards\n    $false,               # MatchSoundsLike\n    $false,               # MatchAllWordForms\n    $true,                # Forward\n    1,                     # Wrap (wdFindContinue)\n    $false,               # Format\n    \"[coordinator]\",      # ReplaceWith\n    2                      # Replace (wdReplaceAll)\n)\n\n# --- Contact person cell --------------------------------------------------\n$contactNeedle = \"[contactname], [contactmail], [contactid], [contactaffiliation], [contactror]\"\n$contactFind = $d.Content.Find\n$contactFind.ClearFormatting()\n$contactFind.Replacement.ClearFormatting()\n$contactFind.Execute(\n    $contactNeedle,\n    $false,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    \"[contact]\",\n    2\n)\n"}
